# Apply updated cryptocurrency price/volume data (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''28.919.12'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.72%  '
$ws.Range("D3").Value = '''1.893.10'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.67%  '
$ws.Range("D4").Value = '''1.001'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.74%  '
$ws.Range("D5").Value = '''325.67'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.25%  '
$ws.Range("D6").Value = '''1.001'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.66%  '
$ws.Range("D7").Value = '''0.4591'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.77%  '
$ws.Range("D8").Value = '''0.3904'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.99%  '
$ws.Range("D9").Value = '''0.07860'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.66%  '
$ws.Range("D10").Value = '''0.9912'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.71%  '
$ws.Range("D11").Value = '''21.90'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.86%  '
$ws.Range("D12").Value = '''1.889.98'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.56%  '
$ws.Range("D13").Value = '''7.041'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.28%  '
$ws.Range("D14").Value = '''5.699'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.20%  '
$ws.Range("D15").Value = '''0.06936'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.48%  '
$ws.Range("D16").Value = '''88.05'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.68%  '
$ws.Range("D17").Value = '''1.002'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.73%  '
$ws.Range("D18").Value = '''0.000009985'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.49%  '
$ws.Range("D19").Value = '''17.02'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.16%  '
$ws.Range("E20").Value = '  -0.46%  '
$ws.Range("D21").Value = '''28.907.52'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.67%  '
$ws.Range("D22").Value = '''5.304'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.00%  '
$ws.Range("D23").Value = '''10.99'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.10%  '
$ws.Range("D24").Value = '''2.086.87'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.19%  '
$ws.Range("D25").Value = '''2.058'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.74%  '
$ws.Range("D26").Value = '''155.84'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.50%  '
$ws.Range("D27").Value = '''19.31'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.19%  '
$ws.Range("D28").Value = '''5.936'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +5.36%  '
$ws.Range("D29").Value = '''1.931'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.31%  '
$ws.Range("D30").Value = '''117.57'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.27%  '
$ws.Range("D31").Value = '''0.09371'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.93%  '
$ws.Range("D32").Value = '''0.9103'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.51%  '
$ws.Range("D33").Value = '''5.297'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.79%  '
$ws.Range("D34").Value = '''1.334'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.66%  '
$ws.Range("D35").Value = '''3.263'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.10%  '
$ws.Range("D36").Value = '''1.189'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.87%  '
$ws.Range("D37").Value = '''0.05775'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.60%  '
$ws.Range("D38").Value = '''0.02075'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.31%  '
$ws.Range("D39").Value = '''1.000'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.69%  '
$ws.Range("D40").Value = '''7.729'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.22%  '
$ws.Range("D41").Value = '''0.5685'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.61%  '
$ws.Range("D42").Value = '''0.1772'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.54%  '
$ws.Range("D43").Value = '''9.780'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.98%  '
$ws.Range("D44").Value = '''2.301'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +9.07%  '
$ws.Range("D45").Value = '''11.89'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.49%  '
$ws.Range("D46").Value = '''0.5361'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.57%  '
$ws.Range("D47").Value = '''0.07040'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.05%  '
$ws.Range("D48").Value = '''1.844'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.17%  '
$ws.Range("B49").Value = 'Quant'
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D49").Value = '''112.74'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.98%  '
$ws.Range("B50").Value = 'MXToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D50").Value = '''2.531'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.97%  '
$ws.Range("D51").Value = '''1.066'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -5.50%  '

Write-Output "Updated cryptos list"
